$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1741
$ws.Range("F4").Value = 147
$ws.Range("F5").Value = 396
$ws.Range("F6").Value = 780
$ws.Range("F7").Value = 225
$ws.Range("F8").Value = 1117
$ws.Range("F11").Value = 860
$ws.Range("F12").Value = 652
$ws.Range("F13").Value = 175
$ws.Range("F14").Value = 501
$ws.Range("F18").Value = 2861
$ws.Range("F19").Value = 2599
$ws.Range("F24").Value = 220
$ws.Range("F25").Value = 18
$ws.Range("F26").Value = 2875
$ws.Range("F32").Value = 1057

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1072
$ws.Range("F5").Value = 1072
$ws.Range("F11").Value = 9
$ws.Range("F18").Value = 976
$ws.Range("F20").Value = 38
$ws.Range("F26").Value = 269
$ws.Range("F27").Value = 3848
$ws.Range("F31").Value = 195
$ws.Range("F34").Value = 145
$ws.Range("F35").Value = 30

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2424
$ws.Range("F6").Value = 1010
$ws.Range("F10").Value = 339
$ws.Range("F11").Value = 92

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2424
$ws.Range("F4").Value = 1741
$ws.Range("F6").Value = 1010
$ws.Range("F8").Value = 339
$ws.Range("F9").Value = 92
$ws.Range("F10").Value = 147
$ws.Range("F11").Value = 396
$ws.Range("F12").Value = 780
$ws.Range("F13").Value = 225
$ws.Range("F15").Value = 1117
$ws.Range("F17").Value = 652
$ws.Range("F18").Value = 1072
$ws.Range("F19").Value = 175
$ws.Range("F20").Value = 501
$ws.Range("F23").Value = 2861
$ws.Range("F24").Value = 2599
$ws.Range("F28").Value = 9
$ws.Range("F29").Value = 220
$ws.Range("F30").Value = 18
$ws.Range("F31").Value = 2877
$ws.Range("F40").Value = 38
$ws.Range("F43").Value = 269
$ws.Range("F44").Value = 1057
$ws.Range("F45").Value = 195
$ws.Range("F48").Value = 145
$ws.Range("F49").Value = 30
